# Update rows 2-6 (existing FAPs sending-cluster rows) and add new rows 7-11 (MuSCs sending-cluster rows)
# reflecting the re-run of the NATMI pipeline with updated TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Adam2/Itgb1 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Adam2"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1251886666666667
$ws.Range("H2").Value = 0.375566
$ws.Range("I2").Value = 0.6155719715657366
$ws.Range("J2").Value = 0.7060466830097307
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 15.264096263758
$ws.Range("R2").Value = 137.376866373822
$ws.Range("S2").Value = 0.1404878145075078
$ws.Range("T2").Value = 0.1708435584343969

# Row 3: FAPs -> Adam2/Itgb1 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Adam2"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1251886666666667
$ws.Range("H3").Value = 0.375566
$ws.Range("I3").Value = 0.6155719715657366
$ws.Range("J3").Value = 0.7060466830097307
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 147.91433
$ws.Range("N3").Value = 443.74299
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 18.51719775359333
$ws.Range("R3").Value = 166.65477978234
$ws.Range("S3").Value = 0.1704287367069571
$ws.Range("T3").Value = 0.2072539311723684

# Row 4: FAPs -> Adam2/Itgb1 -> Inflammatory-Mac
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Adam2"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1251886666666667
$ws.Range("H4").Value = 0.375566
$ws.Range("I4").Value = 0.6155719715657366
$ws.Range("J4").Value = 0.7060466830097307
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 10.45387577088089
$ws.Range("R4").Value = 94.084881937928
$ws.Range("S4").Value = 0.09621546764423161
$ws.Range("T4").Value = 0.1170051148307362

# Row 5: FAPs -> Adam2/Itgb1 -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Adam2"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1251886666666667
$ws.Range("H5").Value = 0.375566
$ws.Range("I5").Value = 0.6155719715657366
$ws.Range("J5").Value = 0.7060466830097307
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 11.40073939575833
$ws.Range("R5").Value = 68.40443637455
$ws.Range("S5").Value = 0.1049302188484368
$ws.Range("T5").Value = 0.08506859729298898

# Row 6: FAPs -> Adam2/Itgb1 -> Resolving-Mac
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Adam2"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1251886666666667
$ws.Range("H6").Value = 0.375566
$ws.Range("I6").Value = 0.6155719715657366
$ws.Range("J6").Value = 0.7060466830097307
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 11.24640273886
$ws.Range("R6").Value = 101.21762464974
$ws.Range("S6").Value = 0.1035097338586033
$ws.Range("T6").Value = 0.1258754812792402

# Row 7: MuSCs -> Adam2/Itgb1 -> ECs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Adam2"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.078181
$ws.Range("H7").Value = 0.156362
$ws.Range("I7").Value = 0.3844280284342634
$ws.Range("J7").Value = 0.2939533169902694
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 9.532510743759
$ws.Range("R7").Value = 57.19506446255399
$ws.Range("S7").Value = 0.08773540064338731
$ws.Range("T7").Value = 0.07112848469754762

# Row 8: MuSCs -> Adam2/Itgb1 -> FAPs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Adam2"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.078181
$ws.Range("H8").Value = 0.156362
$ws.Range("I8").Value = 0.3844280284342634
$ws.Range("J8").Value = 0.2939533169902694
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 147.91433
$ws.Range("N8").Value = 443.74299
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 11.56409023373
$ws.Range("R8").Value = 69.38454140237999
$ws.Range("S8").Value = 0.1064336686320376
$ws.Range("T8").Value = 0.08628746794431304

# Row 9: MuSCs -> Adam2/Itgb1 -> Inflammatory-Mac
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Adam2"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.078181
$ws.Range("H9").Value = 0.156362
$ws.Range("I9").Value = 0.3844280284342634
$ws.Range("J9").Value = 0.2939533169902694
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 6.528502007449333
$ws.Range("R9").Value = 39.171012044696
$ws.Range("S9").Value = 0.06008708037383845
$ws.Range("T9").Value = 0.04871355171970723

# Row 10: MuSCs -> Adam2/Itgb1 -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Adam2"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.078181
$ws.Range("H10").Value = 0.156362
$ws.Range("I10").Value = 0.3844280284342634
$ws.Range("J10").Value = 0.2939533169902694
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 7.119823466712501
$ws.Range("R10").Value = 28.47929386685
$ws.Range("S10").Value = 0.06552948967523393
$ws.Range("T10").Value = 0.03541719966644036

# Row 11: MuSCs -> Adam2/Itgb1 -> Resolving-Mac
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Adam2"
$ws.Range("C11").Value = "Itgb1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.078181
$ws.Range("H11").Value = 0.156362
$ws.Range("I11").Value = 0.3844280284342634
$ws.Range("J11").Value = 0.2939533169902694
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 7.02343938903
$ws.Range("R11").Value = 42.14063633418
$ws.Range("S11").Value = 0.06464238910976602
$ws.Range("T11").Value = 0.05240661296226111
